# "Generate Report for Handoff"
# Status flips from "In Translation" -> "Ready for handoff" and the
# associated "Latest HO/Handback" timestamps advance a bit, on all three
# sheets (Overview, zh-cn, de-de). The two status columns on the Overview
# sheet (and the equivalent Status column on zh-cn / de-de) are also a
# touch wider to fit the new, longer status text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-08-24 12:41:16"
$ws.Columns.Item(5).ColumnWidth = 16.3
$ws.Columns.Item(6).ColumnWidth = 16.3

# --- zh-cn sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-24 12:41:11"
$ws.Columns.Item(3).ColumnWidth = 16.3

# --- de-de sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-24 12:41:16"
$ws.Columns.Item(3).ColumnWidth = 16.3
